$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad) from 46065 -> 46066 for rows 2-13
for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 3).Value = 46066
}

# Rows 6-13 have their A (Beteckning), B (Datum), G (Area ha) values
# cyclically rotated. Capture current values first, then rotate.
$rows = 6..13
$colA = @{}
$colB = @{}
$colG = @{}
foreach ($r in $rows) {
    $colA[$r] = $ws.Cells.Item($r, 1).Value2
    $colB[$r] = $ws.Cells.Item($r, 2).Value2
    $colG[$r] = $ws.Cells.Item($r, 7).Value2
}

# Mapping: new row r gets old data from row srcMap[r]
$srcMap = @{
    6  = 7
    7  = 11
    8  = 12
    9  = 8
    10 = 9
    11 = 13
    12 = 6
    13 = 10
}

foreach ($r in $rows) {
    $src = $srcMap[$r]
    $ws.Cells.Item($r, 1).Value = $colA[$src]
    $ws.Cells.Item($r, 2).Value = $colB[$src]
    $ws.Cells.Item($r, 7).Value = $colG[$src]
}
